# Turn the two-row "frontTest/backTest" sample card deck into a
# practical 4-row Korean/English flash-card deck, and drop the
# special-case "general" alignment on the data rows in favor of "left"
# (removing the old case-sensitivity test pair's quirky formatting).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: frontTest/backTest -> 가게/Store
# (Delete + re-create so the cell picks up the column's default style
# instead of keeping the header-row style it had before.)
$ws.Range("A2:B2").Delete()
$ws.Range("A2").Value = "가게"
$ws.Range("B2").Value = "Store"
$ws.Rows(2).RowHeight = 18.75

# Row 3 (new): 아마도/Maybe
$ws.Range("A3").Value = "아마도"
$ws.Range("B3").Value = "Maybe"
$ws.Rows(3).RowHeight = 18.75

# Row 4 (new): 안녕/Hi
$ws.Range("A4").Value = "안녕"
$ws.Range("B4").Value = "Hi"
$ws.Rows(4).RowHeight = 18.75

# Row 5 (new): 가금/Sometimes, styled like the header row (border + Arial
# + left-aligned), matching the file's existing "header" look.
$ws.Range("A1:B1").Copy($ws.Range("A5:B5"))
$ws.Range("A5").Value = "가금"
$ws.Range("B5").Value = "Sometimes"
$ws.Rows(5).RowHeight = 18.75

# The data rows (2-4) switch from "general" to "left" horizontal
# alignment.
$ws.Range("A2:B4").HorizontalAlignment = -4131
